$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 10098.4
$ws.Range("I19").Value = 15000
$ws.Range("J19").Value = 2746
$ws.Range("K19").Value = 15000
$ws.Range("L19").Value = 2746
$ws.Range("M19").Value = -14825
$ws.Range("N19").Value = -3096
$ws.Range("H33").Value = 450.875
$ws.Range("I33").Value = 472.93332
$ws.Range("J33").Value = 120
$ws.Range("K33").Value = 472.93332
$ws.Range("L33").Value = 120
$ws.Range("M33").Value = -243.93332
$ws.Range("N33").Value = -578
$ws.Range("H62").Value = 2351.4614
$ws.Range("I62").Value = 2367
$ws.Range("K62").Value = 2367
$ws.Range("M62").Value = -1743
$ws.Range("H65").Value = 2351.4614
$ws.Range("I65").Value = 2367
$ws.Range("K65").Value = 11835
$ws.Range("M65").Value = -8715
$ws.Range("H94").Value = 2527.7273
$ws.Range("I94").Value = 2527.7273
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2527.7273
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -2076.7273
$ws.Range("N94").ClearContents()
$ws.Range("H107").Value = 419.73685
$ws.Range("I107").Value = 374.05884
$ws.Range("K107").Value = 374.05884
$ws.Range("M107").Value = 1545.94116
$ws.Range("H113").Value = 1778.2609
$ws.Range("I113").Value = 1545
$ws.Range("J113").Value = 3333.3333
$ws.Range("K113").Value = 1545
$ws.Range("L113").Value = 3333.3333
$ws.Range("M113").Value = 1709
$ws.Range("N113").Value = -9841.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 954.8
$ws.Range("I2").Value = 618.44446
$ws.Range("K2").Value = 618.44446
$ws.Range("M2").Value = -505.44446
$ws.Range("H61").Value = 2081.4
$ws.Range("I61").Value = 1602.4667
$ws.Range("J61").Value = 2799.8
$ws.Range("K61").Value = 1602.4667
$ws.Range("L61").Value = 2799.8
$ws.Range("M61").Value = -1390.4667
$ws.Range("N61").Value = -3223.8
$ws.Range("H116").Value = 954.8
$ws.Range("I116").Value = 618.44446
$ws.Range("K116").Value = 618.44446
$ws.Range("M116").Value = 1675.55554
$ws.Range("H132").Value = 10052.385
$ws.Range("I132").Value = 12998.111
$ws.Range("J132").Value = 3424.5
$ws.Range("K132").Value = 38994.333
$ws.Range("L132").Value = 10273.5
$ws.Range("M132").Value = -36464.333
$ws.Range("N132").Value = -15333.5
$ws.Range("H136").Value = 2081.4
$ws.Range("I136").Value = 1602.4667
$ws.Range("J136").Value = 2799.8
$ws.Range("K136").Value = 4807.4001
$ws.Range("L136").Value = 8399.400000000001
$ws.Range("M136").Value = -2257.4001
$ws.Range("N136").Value = -13499.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 954.8
$ws.Range("I3").Value = 618.44446
$ws.Range("K3").Value = 618.44446
$ws.Range("M3").Value = -504.44446
$ws.Range("H94").Value = 1123.909
$ws.Range("I94").Value = 953.61536
$ws.Range("J94").Value = 1369.8889
$ws.Range("K94").Value = 953.61536
$ws.Range("L94").Value = 1369.8889
$ws.Range("M94").Value = -502.61536
$ws.Range("N94").Value = -2271.8889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H16").Value = 1477.7778
$ws.Range("I16").Value = 1075
$ws.Range("K16").Value = 1075
$ws.Range("M16").Value = -788
$ws.Range("H31").Value = 3848375.2
$ws.Range("I31").Value = 2466
$ws.Range("K31").Value = 2466
$ws.Range("M31").Value = -2171
$ws.Range("H34").Value = 3848375.2
$ws.Range("I34").Value = 2466
$ws.Range("K34").Value = 2466
$ws.Range("M34").Value = -2264
$ws.Range("H113").Value = 1477.7778
$ws.Range("I113").Value = 1075
$ws.Range("K113").Value = 1075
$ws.Range("M113").Value = 1095
$ws.Range("H134").Value = 926
$ws.Range("I134").Value = 872.61536
$ws.Range("K134").Value = 2617.84608
$ws.Range("M134").Value = -82.84608000000026

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 198.33333
$ws.Range("I8").Value = 198.33333
$ws.Range("K8").Value = 594.99999
$ws.Range("M8").Value = -455.99999
$ws.Range("H132").Value = 1067.7142
$ws.Range("J132").Value = 2240
$ws.Range("L132").Value = 20160
$ws.Range("N132").Value = -25220

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 3765000
$ws.Range("I7").Value = 3765000
$ws.Range("K7").Value = 3765000
$ws.Range("M7").Value = -3764888
$ws.Range("H8").Value = 3765000
$ws.Range("I8").Value = 3765000
$ws.Range("K8").Value = 3765000
$ws.Range("M8").Value = -3764861
$ws.Range("H19").Value = 1005.9375
$ws.Range("J19").Value = 1005.9375
$ws.Range("L19").Value = 1005.9375
$ws.Range("N19").Value = -1581.9375
$ws.Range("H24").Value = 600
$ws.Range("I24").Value = 600
$ws.Range("K24").Value = 600
$ws.Range("M24").Value = -427
$ws.Range("H102").Value = 1615.3529
$ws.Range("I102").Value = 1503.8334
$ws.Range("K102").Value = 1503.8334
$ws.Range("M102").Value = 118.1666
$ws.Range("H132").Value = 114901.164
$ws.Range("I132").Value = 157224.84
$ws.Range("J132").Value = 4859.6
$ws.Range("K132").Value = 471674.52
$ws.Range("L132").Value = 14578.8
$ws.Range("M132").Value = -469144.52
$ws.Range("N132").Value = -19638.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 851.6667
$ws.Range("I22").Value = 277.5
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 277.5
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = 17.5
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 851.6667
$ws.Range("I27").Value = 277.5
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 277.5
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -170.5
$ws.Range("N27").Value = -2214
$ws.Range("H122").Value = 2413.2173
$ws.Range("I122").Value = 3700.5715
$ws.Range("J122").Value = 1850
$ws.Range("K122").Value = 11101.7145
$ws.Range("L122").Value = 5550
$ws.Range("M122").Value = -8651.7145
$ws.Range("N122").Value = -10450
$ws.Range("H132").Value = 16028.782
$ws.Range("I132").Value = 25016.924
$ws.Range("J132").Value = 4344.2
$ws.Range("K132").Value = 75050.772
$ws.Range("L132").Value = 13032.6
$ws.Range("M132").Value = -72520.772
$ws.Range("N132").Value = -18092.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 1498
$ws.Range("J22").Value = 1498
$ws.Range("L22").Value = 1498
$ws.Range("N22").Value = -2084
$ws.Range("H42").Value = 51699.332
$ws.Range("I42").Value = 15000
$ws.Range("J42").Value = 70049
$ws.Range("K42").Value = 15000
$ws.Range("L42").Value = 70049
$ws.Range("M42").Value = -14622
$ws.Range("N42").Value = -70805
$ws.Range("H126").Value = 1266.5264
$ws.Range("I126").Value = 1716.4
$ws.Range("J126").Value = 766.6667
$ws.Range("K126").Value = 5149.200000000001
$ws.Range("L126").Value = 2300.0001
$ws.Range("M126").Value = -2679.200000000001
$ws.Range("N126").Value = -7240.0001
